# Refresh the cryptos worksheet with newly scraped price/volume figures.
# (GitHub Actions data refresh - Sun Feb 25 23:19:28 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.718.48"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "3.100.20"
$ws.Range("E3").Value = "  +3.77%  "
$ws.Range("D4").Value = "'" + "1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'" + "389.53"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").Value = "'" + "103.30"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("D7").Value = "'" + "0.543"
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'" + "0.591"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").Value = "'" + "37.30"
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("D11").Value = "'" + "0.137"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "3.598.58"
$ws.Range("E13").Value = "  +3.74%  "
$ws.Range("D14").Value = "'" + "18.67"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "'" + "7.85"
$ws.Range("D16").Value = "3.108.71"
$ws.Range("E16").Value = "  +4.21%  "
$ws.Range("D17").Value = "'" + "0.987"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").Value = "'" + "10.83"
$ws.Range("E18").Value = "  -3.90%  "
$ws.Range("D19").Value = "51.842.60"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").Value = "'" + "3.20"
$ws.Range("E20").Value = "  +3.87%  "
$ws.Range("D21").Value = "'" + "12.45"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "0.0₃0967"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'" + "70.04"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").Value = "'" + "268.10"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E26").Value = "  +1.75%  "
$ws.Range("D27").Value = "'" + "27.13"
$ws.Range("E27").Value = "  +3.76%  "
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "'" + "7.14"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").Value = "'" + "0.109"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").Value = "'" + "10.36"
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("D33").Value = "'" + "35.39"
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("E34").Value = "  +0.78%  "
$ws.Range("D35").Value = "'" + "50.28"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").Value = "'" + "0.0448"
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "'" + "3.40"
$ws.Range("E38").Value = "  +2.86%  "
$ws.Range("D39").Value = "'" + "0.289"
$ws.Range("E39").Value = "  +6.06%  "
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("D41").Value = "'" + "2.61"
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("D42").Value = "'" + "16.85"
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("D43").Value = "'" + "129.05"
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").Value = "'" + "3.71"
$ws.Range("E45").Value = "  -3.87%  "
$ws.Range("D46").Value = "'" + "22.20"
$ws.Range("E46").Value = "  +3.80%  "
$ws.Range("E47").Value = "  +6.28%  "
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("D49").Value = "2.047.53"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").Value = "3.409.81"
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'" + "0.206"
$ws.Range("E51").Value = "  +4.64%  "
